# Function Progress Checker.xlsx — add an "Edge Tested" column (G) next to the
# existing "Tested" column (F), which is renamed "Primary Tested". A few rows'
# F-column entries are adjusted to reflect the new primary/edge-tested split.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column G: width + header -------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 12.736979166666666

# Give G2 the same header formatting as the other header cells (row 2), then
# set its text. Rename F2 from "Tested" to "Primary Tested" and add the new
# "Edge Tested" header in G2.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Primary Tested"
$ws.Range("G2").Value = "Edge Tested"

# --- Row 3 (ADD) and Row 6 (ADDU): no longer marked tested in column F ------------
$ws.Range("F3").ClearContents()
$ws.Range("F6").ClearContents()

# --- Row 22 (JR): note the edge case for the existing "tested" mark ---------------
$ws.Range("F22").Value = "d (jr = 0 only)"

# --- Row 27 (LUI) and Row 37 (OR): now also marked as tested in column F ----------
$ws.Range("F3").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F37").PasteSpecial(-4122)
$ws.Range("F27").Value = "d"
$ws.Range("F37").Value = "d"

# --- Final selection, matching where editing left off ------------------------------
[void]$ws.Range("F24").Select()
